$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" summary text (A1) with the new conversion rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.68 = 18514.36 pesos`n✅ 18514.36 pesos = 4.65 = 910.85 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Update the "tasas" sheet rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 213.888
$wsTasas.Range("O10").Value = 3960
$wsTasas.Range("N12").Value = 3983.99
$wsTasas.Range("O12").Value = 196
